$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 219 (shifts existing rows 219:315 down to 220:316)
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new data record
$ws.Range("A219").Value = 5
$ws.Range("B219").Value = "Macroferia Regional de Talca"
$ws.Range("C219").Value = "Maule"
$ws.Range("D219").Value = 44784
$ws.Range("E219").Value = 7
$ws.Range("F219").Value = "Fruta"
$ws.Range("G219").Value = 100101
$ws.Range("H219").Value = "Berries"
$ws.Range("I219").Value = 100101007
$ws.Range("J219").Value = "Kiwi"
$ws.Range("K219").Value = "Hayward"
$ws.Range("L219").Value = "Primera"
$ws.Range("M219").Value = 300
$ws.Range("N219").Value = 6000
$ws.Range("O219").Value = 6000
$ws.Range("P219").Value = 6000
$ws.Range("Q219").Value = "$/bandeja 18 kilos"
$ws.Range("R219").Value = "Provincia de Curicó"
$ws.Range("S219").Value = 333
$ws.Range("T219").Value = 18
